# Update the "dSF" column (F) values to match the repulled/pushed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 1
    16 = -1
    17 = 4
    20 = -4
    23 = -2
    27 = -3
    28 = 0
    30 = -2
    32 = -7
    33 = 2
    35 = -4
    42 = -2
    45 = 0
    46 = -9
    51 = 1
    52 = 2
    55 = 2
    60 = -4
    63 = -2
    70 = -1
    73 = 0
    76 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
